$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 33: "good words" method assignment entry, mirroring the existing
# "bad words" row (row 30) structure. Shared-string insertion order below
# matches column order E, B, C, D, F so new <si> entries land in the same
# sequence as the target workbook.
$ws.Range("E33").Value = "def getQuotesContainingGoodWords()"
$ws.Range("B33").Value = "TBD"
$ws.Range("C33").Value = "TBD"
$ws.Range("D33").Value = "Return the list of Quotes which contains some good words"
$ws.Range("F33").Value = "crypto.py"
$ws.Range("A33").Value = 32

# Reflect the author's scroll/selection after adding the row: the view was
# scrolled down so the new row is visible, with E35 left as the active cell.
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E35").Select()
